$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.730.10"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "3.008.90"
$ws.Range("E3").Value = "  +2.93%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.48%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.08%  "

$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0849"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "3.487.93"
$ws.Range("E14").Value = "  +2.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.14%  "

$ws.Range("D16").Value = "3.019.96"
$ws.Range("E16").Value = "  +3.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.981"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").Value = "51.784.43"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("E22").Value = "  +1.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.18%  "

$ws.Range("E26").Value = "  +19.37%  "

$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "26.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").Value = "  -3.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("E35").Value = "  -3.75%  "

$ws.Range("E36").Value = "  +5.13%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.74%  "

$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("E42").Value = "  +2.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("E45").Value = "  -2.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.277"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.68%  "

$ws.Range("D48").Value = "2.060.51"
$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0355"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.35%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.69%  "
